$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Arduino-side Protocols:" -> "Arduino to Raspberry Pi Protocols:"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Arduino-side Protocols", $false, $false, $false, $false, $false, $true, 1, $false, "Arduino to Raspberry Pi Protocols", 2)

# ---------------------------------------------------------------------------
# 2) Underscore -> semicolon separators in the "String summary" paragraph
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("S_L1_A_L2_R", $false, $false, $false, $false, $false, $true, 1, $false, "S;L1;A;L2;R", 2)
$d.Content.Find.Execute("S_1_1_2_10", $false, $false, $false, $false, $false, $true, 1, $false, "S;1;1;2;10", 2)
$d.Content.Find.Execute("S_5_01011_2_10", $false, $false, $false, $false, $false, $true, 1, $false, "S;5;01011;2;10", 2)

# ---------------------------------------------------------------------------
# 3) "Raspberry Pi-side Protocols:" -> "Raspberry Pi to Server Protocols:"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Raspberry Pi-side Protocols", $false, $false, $false, $false, $false, $true, 1, $false, "Raspberry Pi to Server Protocols", 2)

# ---------------------------------------------------------------------------
# 4) "Server-side Protocols:" -> "Server to Raspberry Pi Protocols:"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Server-side Protocols", $false, $false, $false, $false, $false, $true, 1, $false, "Server to Raspberry Pi Protocols", 2)

# ---------------------------------------------------------------------------
# 5) "u " update." -> "u " update (seats)."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("u " + [char]8211 + " update.", $false, $false, $false, $false, $false, $true, 1, $false, "u " + [char]8211 + " update (seats).", 2)

# ---------------------------------------------------------------------------
# 6) "c " close." -> "c " close (socket)."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("c " + [char]8211 + " close.", $false, $false, $false, $false, $false, $true, 1, $false, "c " + [char]8211 + " close (socket).", 2)

# ---------------------------------------------------------------------------
# 7) Extend the "received" sentence with the extra clause about retrying
# ---------------------------------------------------------------------------
$oldAck = "for ack )"
$newAck = "for ack. If RPI doesn" + [char]8217 + "t get " + [char]8216 + "r" + [char]8217 + ", doesn" + [char]8217 + "t continue  )"
$d.Content.Find.Execute($oldAck, $false, $false, $false, $false, $false, $true, 1, $false, $newAck, 2)

# ---------------------------------------------------------------------------
# 8) Move the "_GoBack" bookmark from the end of the "received" paragraph to
#    right after "(seats)." in the "update" line (the last edited spot).
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()
$bmRange = $d.Content
$bmRange.Find.Execute("update (seats).")
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# 9) Turn the trailing empty paragraph into the new "Client to Server
#    Protocols" section.
# ---------------------------------------------------------------------------
$dash = [char]8211
$lineBreak = [char]11

$newPara = $d.Paragraphs.Last
$newRng = $newPara.Range
$newRng.Collapse(1)  # wdCollapseStart
$insertStart = $newRng.Start

$headerText = "Client to Server Protocols:" + "$lineBreak$lineBreak"
$bodyText = "c " + $dash + " close (socket)." + "$lineBreak" + "g " + $dash + " get (seats)."
$newRng.InsertAfter($headerText + $bodyText)

$headerRange = $d.Range($insertStart, $insertStart + $headerText.Length)
$headerRange.Font.Underline = 1

Write-Output $d.Content.Text
